$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.735.96"
$ws.Range("E2").Value = "  -0.16%  "

$ws.Range("D3").Value = "1.534.10"
$ws.Range("E3").Value = "  -1.98%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("E5").Value = "  -0.32%  "

$ws.Range("E6").Value = "  -1.01%  "

$ws.Range("E7").Value = "  -0.06%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "21.31"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  -2.98%  "

$ws.Range("E9").Value = "  -1.15%  "

$ws.Range("E10").Value = "  -0.70%  "

$ws.Range("E11").Value = "  -0.81%  "

$ws.Range("E12").Value = "  -1.91%  "

$ws.Range("D13").Value = "1.547.38"
$ws.Range("E13").Value = "  -1.25%  "

$ws.Range("E14").Value = "  -1.58%  "

$ws.Range("E15").Value = "  -1.43%  "

$ws.Range("D16").Value = "26.722.20"
$ws.Range("E16").Value = "  -0.20%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "61.00"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -0.89%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "212.70"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  -0.70%  "

$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "7.22"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -1.90%  "

$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "0.0₃0681"
$ws.Range("E20").Value = "  +0.81%  "

$ws.Range("E21").Value = "  +0.01%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "3.99"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -2.34%  "

$ws.Range("E23").Value = "  -3.08%  "

$ws.Range("E24").Value = "  -2.93%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "151.90"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -0.63%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "6.58"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  -2.42%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "14.79"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -0.52%  "

$ws.Range("E28").Value = "  -0.02%  "

$ws.Range("E29").Value = "  -1.32%  "

$ws.Range("E30").Value = "  -1.33%  "

$ws.Range("E31").Value = "  -1.99%  "

$ws.Range("E32").Value = "  +2.60%  "

$ws.Range("D33").Value = "1.363.53"
$ws.Range("E33").Value = "  -1.61%  "

$ws.Range("E34").Value = "  -0.25%  "

$ws.Range("E35").Value = "  -2.87%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "0.945"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  +2.55%  "

$ws.Range("E37").Value = "  -0.23%  "

$ws.Range("E38").Value = "  -0.09%  "

$ws.Range("E39").Value = "  -0.05%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.799"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -1.83%  "

$ws.Range("E41").Value = "  +6.08%  "

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "0.998"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +0.58%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "2.19"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  +0.61%  "

$ws.Range("E45").Value = "  -1.22%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "62.43"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -1.44%  "

$ws.Range("D47").Value = "1.666.29"
$ws.Range("E47").Value = "  -1.82%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "85.34"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -0.31%  "

$ws.Range("E49").Value = "  +2.13%  "

$ws.Range("D50").Value = "0.0₇0973"
$ws.Range("E50").Value = "  -1.16%  "

$ws.Range("E51").Value = "  -0.73%  "
